{"js": "// Add a new \"Mini Project 2\" status row (Sr No. 10) to the Core Java\n// assignment-status table (the 2nd table in the document body), then\n// relocate the trailing \"_GoBack\" bookmark from the paragraph that used\n// to sit right after the table into the new row's last cell, matching\n// where Word leaves the cursor after the edit.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// The Core Java table is the second table in the document.\nconst table = tables.items[1];\n\n// Append the new row; Word inherits the run/paragraph formatting of the\n// preceding row (bold, accent1 color, sz 28/36, single underline, en-US).\ntable.addRows(\"End\", 1, [\n  [\"10\", \"Mini Project 2\", \"21/01/2021\", \"Done till sprint 2\", \"24/01/2021\", \"done\"]\n]);\nawait context.sync();\n\n// Grab the newly added (last) row/cell.\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst lastRow = rows.items[rows.items.length - 1];\nconst cells = lastRow.cells;\ncells.load(\"items\");\nawait context.sync();\n\nconst lastCell = cells.items[cells.items.length - 1];\n\n// Move the \"_GoBack\" bookmark into the start of that last cell.\ncontext.document.deleteBookmark(\"_GoBack\");\nconst insertionPoint = lastCell.body.getRange(\"Start\");\ninsertionPoint.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Add a new \"Mini Project 2\" status row (Sr No. 10) to the Core Java\n# assignment-status table (the 2nd table in the document), then relocate\n# the trailing \"_GoBack\" bookmark from the paragraph that used to sit\n# right after the table into the new row's last cell, matching where\n# Word leaves the cursor after the edit.\n\n$d = $word.ActiveDocument\n\n# The Core Java table is the second table in the document.\n$table = $d.Tables.Item(2)\n\n# Append a new row; Word inherits the row/paragraph/run formatting\n# (bold, accent1 color, sz 28/36, single underline, en-US) from the\n# preceding row automatically.\n$newRow = $table.Rows.Add()\n\n$newRow.Cells.Item(1).Range.Text = \"10\"\n$newRow.Cells.Item(2).Range.Text = \"Mini Project 2\"\n$newRow.Cells.Item(3).Range.Text = \"21/01/2021\"\n$newRow.Cells.Item(4).Range.Text = \"Done till sprint 2\"\n$newRow.Cells.Item(5).Range.Text = \"24/01/2021\"\n$newRow.Cells.Item(6).Range.Text = \"done\"\n\n# Move the \"_GoBack\" bookmark to the start of the new last cell.\n# (Re-fetch the position through $d.Range(...) rather than reusing the\n# cell's own Range object for the Add call.)\n$lastCell = $newRow.Cells.Item(6)\n$startPos = $lastCell.Range.Start\n$bookmarkRange = $d.Range($startPos, $startPos)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n"}
